$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (theta_se) - update standard errors
$ws.Range("C4").Value = "(0.37)"
$ws.Range("D4").Value = "(0.18)"
$ws.Range("E4").Value = "(0.2)"
$ws.Range("F4").Value = "(0.26)"
$ws.Range("G4").Value = "(0.28)"

# Row 6 (lambda_se) - update standard errors
$ws.Range("C6").Value = "(0.22)"
$ws.Range("D6").Value = "(0.14)"
$ws.Range("E6").Value = "(0.17)"
$ws.Range("F6").Value = "(0.19)"
$ws.Range("G6").Value = "(0.19)"
